$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (shared string) renames ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C numeric updates (rows 2-137) ---
$ws.Range("C2").Value = 2934.187009790061
$ws.Range("C3").Value = 2870.311589353206
$ws.Range("C4").Value = 1873.394108966653
$ws.Range("C5").Value = 8947.741473873051
$ws.Range("C6").Value = 1460.056109840828
$ws.Range("C7").Value = 1909.084588129339
$ws.Range("C8").Value = 6128.19547247793
$ws.Range("C9").Value = 4547.50930098406
$ws.Range("C10").Value = 3972.630273980753
$ws.Range("C11").Value = 4729.735976516416
$ws.Range("C12").Value = 748.2960604568028
$ws.Range("C13").Value = 1268.249210347625
$ws.Range("C14").Value = 567.9059336271471
$ws.Range("C15").Value = 993.3829437244538
$ws.Range("C16").Value = 1250.795760575873
$ws.Range("C17").Value = 478.6685897045245
$ws.Range("C18").Value = 471.181692645893
$ws.Range("C19").Value = 17288.8595992193
$ws.Range("C20").Value = 341.5541149051794
$ws.Range("C21").Value = 951.6879611168786
$ws.Range("C22").Value = 815.8736791314819
$ws.Range("C23").Value = 2983.242707849043
$ws.Range("C24").Value = 2898.942214704482
$ws.Range("C25").Value = 1904.346464968814
$ws.Range("C26").Value = 9271.398233246389
$ws.Range("C27").Value = 1503.870423231357
$ws.Range("C28").Value = 1955.461557360978
$ws.Range("C29").Value = 6336.709213679884
$ws.Range("C30").Value = 4633.590358399045
$ws.Range("C31").Value = 4355.934938677345
$ws.Range("C32").Value = 5082.354756663512
$ws.Range("C33").Value = 781.1535935570469
$ws.Range("C34").Value = 1357.563719132622
$ws.Range("C35").Value = 592.4010974509293
$ws.Range("C36").Value = 987.4097230439231
$ws.Range("C37").Value = 1317.890706178356
$ws.Range("C38").Value = 487.7306818514292
$ws.Range("C39").Value = 492.3430015592067
$ws.Range("C40").Value = 17610.30663334184
$ws.Range("C41").Value = 369.2024078290272
$ws.Range("C42").Value = 982.980837581714
$ws.Range("C43").Value = 864.5379000312432
$ws.Range("C44").Value = 3083.80337578809
$ws.Range("C45").Value = 2965.153206179127
$ws.Range("C46").Value = 691.8942672110555
$ws.Range("C47").Value = 1939.33862702996
$ws.Range("C48").Value = 9477.887185090232
$ws.Range("C49").Value = 1577.487171555845
$ws.Range("C50").Value = 2024.117324382548
$ws.Range("C51").Value = 6711.616186806423
$ws.Range("C52").Value = 4921.848409120176
$ws.Range("C53").Value = 4479.398934239905
$ws.Range("C54").Value = 5360.226632400601
$ws.Range("C55").Value = 822.1883388417289
$ws.Range("C56").Value = 1410.426304742003
$ws.Range("C57").Value = 612.032557723897
$ws.Range("C58").Value = 992.8781394745556
$ws.Range("C59").Value = 1385.890384668919
$ws.Range("C60").Value = 2094.024217383061
$ws.Range("C61").Value = 6911.59200404802
$ws.Range("C62").Value = 5122.180090208862
$ws.Range("C63").Value = 3156.723844635973
$ws.Range("C64").Value = 2999.422762626143
$ws.Range("C65").Value = 1982.009737844954
$ws.Range("C66").Value = 9690.869064532331
$ws.Range("C67").Value = 1657.651524528445
$ws.Range("C68").Value = 4394.543881413723
$ws.Range("C69").Value = 2201.396847776877
$ws.Range("C70").Value = 7200.731056811853
$ws.Range("C71").Value = 5295.682695961288
$ws.Range("C72").Value = 3212.740625904757
$ws.Range("C73").Value = 3056.152683606517
$ws.Range("C74").Value = 2000.792448761861
$ws.Range("C75").Value = 9693.722968944676
$ws.Range("C76").Value = 1716.389195271215
$ws.Range("C77").Value = 4699.493713911862
$ws.Range("C78").Value = 951.3148210424945
$ws.Range("C79").Value = 2286.013198234259
$ws.Range("C80").Value = 7449.08671983612
$ws.Range("C81").Value = 5412.131646018807
$ws.Range("C82").Value = 3252.634165082374
$ws.Range("C83").Value = 449.4203771491282
$ws.Range("C84").Value = 3137.260298393558
$ws.Range("C85").Value = 730.3063521039821
$ws.Range("C86").Value = 2025.814194788851
$ws.Range("C87").Value = 1640.18070024053
$ws.Range("C88").Value = 1060.095015975378
$ws.Range("C89").Value = 507.537974993908
$ws.Range("C90").Value = 9839.050190896
$ws.Range("C91").Value = 558.2093442539386
$ws.Range("C92").Value = 711.3043470146426
$ws.Range("C93").Value = 1775.027517189621
$ws.Range("C94").Value = 1054.227974015008
$ws.Range("C95").Value = 4861.287098802361
$ws.Range("C96").Value = 5996.49696468919
$ws.Range("C97").Value = 886.4370030633224
$ws.Range("C98").Value = 1579.189101937001
$ws.Range("C99").Value = 1002.388731936373
$ws.Range("C100").Value = 2361.056581219794
$ws.Range("C101").Value = 7580.275568826287
$ws.Range("C102").Value = 5330.539154475424
$ws.Range("C103").Value = 3314.741082534716
$ws.Range("C104").Value = 482.6390663355013
$ws.Range("C105").Value = 3210.869677115934
$ws.Range("C106").Value = 729.1196658666737
$ws.Range("C107").Value = 2067.29003376698
$ws.Range("C108").Value = 1751.664428859304
$ws.Range("C109").Value = 1093.134170274031
$ws.Range("C110").Value = 507.5484050163182
$ws.Range("C111").Value = 10037.20149040966
$ws.Range("C112").Value = 579.0880693780265
$ws.Range("C113").Value = 731.9993357350996
$ws.Range("C114").Value = 1836.014008604312
$ws.Range("C115").Value = 1081.294365994858
$ws.Range("C116").Value = 4944.191641077407
$ws.Range("C117").Value = 6114.227214287786
$ws.Range("C118").Value = 900.3889853519216
$ws.Range("C119").Value = 1667.171891046301
$ws.Range("C120").Value = 2425.561644739583
$ws.Range("C121").Value = 7633.969039669125
$ws.Range("C122").Value = 5176.058803160127
$ws.Range("C123").Value = 3382.563653843273
$ws.Range("C124").Value = 514.0573067519859
$ws.Range("C125").Value = 3242.636921959078
$ws.Range("C126").Value = 729.8559996981501
$ws.Range("C127").Value = 2111.193164269742
$ws.Range("C128").Value = 1875.732161108182
$ws.Range("C129").Value = 1129.713195979213
$ws.Range("C130").Value = 506.2496613373833
$ws.Range("C131").Value = 10205.79575322194
$ws.Range("C132").Value = 584.2111078769213
$ws.Range("C133").Value = 729.6614300490079
$ws.Range("C134").Value = 1895.214690888655
$ws.Range("C135").Value = 5089.61202008711
$ws.Range("C136").Value = 6262.368904654469
$ws.Range("C137").Value = 909.5979669529498

# --- AL column flag updates ---
$ws.Range("AL46").Value = 1
$ws.Range("AL85").Value = 1
$ws.Range("AL106").Value = 1
$ws.Range("AL126").Value = 1
